$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Adam10", "Epha3", "ECs", 3, 1, 68.44130066666666, 205.323902, 0.4535705444232658, 0.4535705444232658, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.6741239977997777, 6.067115980197999, 0.00009134349770196041, 0.00009134349770196041),
  @("ECs", "Adam10", "Epha3", "FAPs", 3, 1, 68.44130066666666, 205.323902, 0.4535705444232658, 0.4535705444232658, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 3165.47412563569, 28489.26713072121, 0.4289203165357394, 0.4289203165357394),
  @("ECs", "Adam10", "Epha3", "MuSCs", 3, 1, 68.44130066666666, 205.323902, 0.4535705444232658, 0.4535705444232658, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 178.4878170571642, 1606.390353514478, 0.02418501872118694, 0.02418501872118695),
  @("ECs", "Adam10", "Epha3", "Resolving-Mac", 3, 1, 68.44130066666666, 205.323902, 0.4535705444232658, 0.4535705444232658, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 2.759165408842889, 24.832488679586, 0.0003738656686374564, 0.0003738656686374564),
  @("FAPs", "Adam10", "Epha3", "ECs", 3, 1, 22.11880033333333, 66.35640099999999, 0.1465845361127928, 0.1465845361127928, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.2178628103498888, 1.960765293149, 0.00002952031255598222, 0.00002952031255598222),
  @("FAPs", "Adam10", "Epha3", "FAPs", 3, 1, 22.11880033333333, 66.35640099999999, 0.1465845361127928, 0.1465845361127928, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 1023.015189122045, 9207.136702098403, 0.1386180967917338, 0.1386180967917338),
  @("FAPs", "Adam10", "Epha3", "MuSCs", 3, 1, 22.11880033333333, 66.35640099999999, 0.1465845361127928, 0.1465845361127928, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 57.68353828703209, 519.151844583289, 0.007816093425185286, 0.007816093425185286),
  @("FAPs", "Adam10", "Epha3", "Resolving-Mac", 3, 1, 22.11880033333333, 66.35640099999999, 0.1465845361127928, 0.1465845361127928, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 0.8917046895714443, 8.025342206142998, 0.0001208255833178164, 0.0001208255833178164),
  @("MuSCs", "Adam10", "Epha3", "ECs", 3, 1, 28.16264866666667, 84.48794600000001, 0.1866380060535937, 0.1866380060535937, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.2773927018171111, 2.496534316354, 0.0000375865860044602, 0.0000375865860044602),
  @("MuSCs", "Adam10", "Epha3", "FAPs", 3, 1, 28.16264866666667, 84.48794600000001, 0.1866380060535937, 0.1866380060535937, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 1302.548823522287, 11722.93941170058, 0.1764947781957429, 0.1764947781957429),
  @("MuSCs", "Adam10", "Epha3", "MuSCs", 3, 1, 28.16264866666667, 84.48794600000001, 0.1866380060535937, 0.1866380060535937, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 73.44526819475489, 661.007413752794, 0.009951800719843282, 0.009951800719843284),
  @("MuSCs", "Adam10", "Epha3", "Resolving-Mac", 3, 1, 28.16264866666667, 84.48794600000001, 0.1866380060535937, 0.1866380060535937, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 1.135358405897556, 10.218225653078, 0.0001538405520030265, 0.0001538405520030265),
  @("Resolving-Mac", "Adam10", "Epha3", "ECs", 3, 1, 32.17175066666667, 96.515252, 0.2132069134103476, 0.2132069134103476, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.3168810201497778, 2.851929181348, 0.00004293723533106307, 0.00004293723533106307),
  @("Resolving-Mac", "Adam10", "Epha3", "FAPs", 3, 1, 32.17175066666667, 96.515252, 0.2132069134103476, 0.2132069134103476, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 1487.97353819629, 13391.76184376661, 0.2016197434157795, 0.2016197434157795),
  @("Resolving-Mac", "Adam10", "Epha3", "MuSCs", 3, 1, 32.17175066666667, 96.515252, 0.2132069134103476, 0.2132069134103476, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 83.90059060051422, 755.1053154046281, 0.01136849219093876, 0.01136849219093876),
  @("Resolving-Mac", "Adam10", "Epha3", "Resolving-Mac", 3, 1, 32.17175066666667, 96.515252, 0.2132069134103476, 0.2132069134103476, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 1.296982680292889, 11.672844122636, 0.0001757405682982423, 0.0001757405682982423)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $rowData = $rows[$i]
  for ($j = 0; $j -lt $rowData.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $rowData[$j]
  }
}

Write-Host "Updated $($rows.Count) rows"
